$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 377.34784
$ws.Range("I6").Value = 377.34784
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1132.04352
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1020.04352
$ws.Range("H55").Value = 240
$ws.Range("I55").Value = 266.66666
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 266.66666
$ws.Range("L55").Value = 200
$ws.Range("M55").Value = -52.66665999999998
$ws.Range("N55").Value = -628
$ws.Range("H64").Value = 2962964.5
$ws.Range("J64").Value = 5675
$ws.Range("L64").Value = 5675
$ws.Range("N64").Value = -6171
$ws.Range("H67").Value = 2962964.5
$ws.Range("J67").Value = 5675
$ws.Range("L67").Value = 5675
$ws.Range("N67").Value = -7391
$ws.Range("H76").Value = 4350.625
$ws.Range("I76").Value = 2944.2856
$ws.Range("J76").Value = 5444.4443
$ws.Range("K76").Value = 2944.2856
$ws.Range("L76").Value = 5444.4443
$ws.Range("M76").Value = -2629.2856
$ws.Range("N76").Value = -6074.4443
$ws.Range("H79").Value = 4350.625
$ws.Range("I79").Value = 2944.2856
$ws.Range("J79").Value = 5444.4443
$ws.Range("K79").Value = 2944.2856
$ws.Range("L79").Value = 5444.4443
$ws.Range("M79").Value = -1852.2856
$ws.Range("N79").Value = -7628.4443
$ws.Range("H103").Value = 6010380.5
$ws.Range("I103").Value = 15025201
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 45075603
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = -45075017
$ws.Range("N103").Value = -2672
$ws.Range("H129").Value = 2117.74
$ws.Range("I129").Value = 458.7
$ws.Range("J129").Value = 2532.5
$ws.Range("K129").Value = 1376.1
$ws.Range("L129").Value = 7597.5
$ws.Range("M129").Value = 3623.9
$ws.Range("N129").Value = -17597.5
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4937.057
$ws.Range("I122").Value = 5955.56
$ws.Range("K122").Value = 17866.68
$ws.Range("M122").Value = -15416.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 33071.332
$ws.Range("J76").Value = 33071.332
$ws.Range("L76").Value = 33071.332
$ws.Range("N76").Value = -33701.332
$ws.Range("H79").Value = 33071.332
$ws.Range("J79").Value = 33071.332
$ws.Range("L79").Value = 33071.332
$ws.Range("N79").Value = -35255.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 30092
$ws.Range("J54").Value = 30092
$ws.Range("L54").Value = 30092
$ws.Range("N54").Value = -31408
$ws.Range("H110").Value = 40702
$ws.Range("J110").Value = 40702
$ws.Range("L110").Value = 40702
$ws.Range("N110").Value = -48882
$ws.Range("H132").Value = 8198373
$ws.Range("I132").Value = 9805332
$ws.Range("K132").Value = 29415996
$ws.Range("M132").Value = -29413466

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1369.2858
$ws.Range("I7").Value = 300
$ws.Range("J7").Value = 1797
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 5391
$ws.Range("M7").Value = -788
$ws.Range("N7").Value = -5615
$ws.Range("H68").Value = 1002.9286
$ws.Range("I68").Value = 1533.3334
$ws.Range("J68").Value = 858.2727
$ws.Range("K68").Value = 4600.0002
$ws.Range("L68").Value = 2574.8181
$ws.Range("M68").Value = -3789.0002
$ws.Range("N68").Value = -4196.8181
$ws.Range("H71").Value = 1002.9286
$ws.Range("I71").Value = 1533.3334
$ws.Range("J71").Value = 858.2727
$ws.Range("K71").Value = 13800.0006
$ws.Range("L71").Value = 7724.454299999999
$ws.Range("M71").Value = -9744.000599999999
$ws.Range("N71").Value = -15836.4543
$ws.Range("H80").Value = 2867.6667
$ws.Range("J80").Value = 2867.6667
$ws.Range("L80").Value = 8603.000100000001
$ws.Range("N80").Value = -10475.0001
$ws.Range("H83").Value = 2867.6667
$ws.Range("J83").Value = 2867.6667
$ws.Range("L83").Value = 25809.0003
$ws.Range("N83").Value = -35169.0003
$ws.Range("H92").Value = 825.75
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 901
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 2703
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = -5199
$ws.Range("H113").Value = 714.45
$ws.Range("I113").Value = 437.47827
$ws.Range("J113").Value = 1089.1765
$ws.Range("K113").Value = 1312.43481
$ws.Range("L113").Value = 3267.5295
$ws.Range("M113").Value = 857.56519
$ws.Range("N113").Value = -7607.529500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7578678.5
$ws.Range("I80").Value = 16669147
$ws.Range("J80").Value = 3287.7083
$ws.Range("K80").Value = 16669147
$ws.Range("L80").Value = 3287.7083
$ws.Range("M80").Value = -16668149
$ws.Range("N80").Value = -5283.7083
$ws.Range("H83").Value = 7578678.5
$ws.Range("I83").Value = 16669147
$ws.Range("J83").Value = 3287.7083
$ws.Range("K83").Value = 83345735
$ws.Range("L83").Value = 16438.5415
$ws.Range("M83").Value = -83340743
$ws.Range("N83").Value = -26422.5415
$ws.Range("H97").Value = 1445.909
$ws.Range("I97").Value = 1361.75
$ws.Range("J97").Value = 1670.3334
$ws.Range("K97").Value = 1361.75
$ws.Range("L97").Value = 1670.3334
$ws.Range("M97").Value = -865.75
$ws.Range("N97").Value = -2662.3334
$ws.Range("H132").Value = 3826.5
$ws.Range("I132").Value = 2712
$ws.Range("K132").Value = 8136
$ws.Range("M132").Value = -5606

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8239.111000000001
$ws.Range("I7").Value = 21001.334
$ws.Range("J7").Value = 5686.6665
$ws.Range("K7").Value = 21001.334
$ws.Range("L7").Value = 5686.6665
$ws.Range("M7").Value = -20889.334
$ws.Range("N7").Value = -5910.6665
$ws.Range("H40").Value = 6018
$ws.Range("I40").Value = 11900
$ws.Range("J40").Value = 3755.6924
$ws.Range("K40").Value = 11900
$ws.Range("L40").Value = 3755.6924
$ws.Range("M40").Value = -11764
$ws.Range("N40").Value = -4027.6924
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H68").Value = 1901
$ws.Range("I68").Value = 1934.6666
$ws.Range("J68").Value = 1800
$ws.Range("K68").Value = 1934.6666
$ws.Range("L68").Value = 1800
$ws.Range("M68").Value = -1185.6666
$ws.Range("N68").Value = -3298
$ws.Range("H71").Value = 1901
$ws.Range("I71").Value = 1934.6666
$ws.Range("J71").Value = 1800
$ws.Range("K71").Value = 9673.333000000001
$ws.Range("L71").Value = 9000
$ws.Range("M71").Value = -5929.333000000001
$ws.Range("N71").Value = -16488
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20676
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22340
$ws.Range("H82").Value = 2227.2727
$ws.Range("I82").Value = 1750
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 1750
$ws.Range("L82").Value = 2800
$ws.Range("M82").Value = -1389
$ws.Range("N82").Value = -3522
$ws.Range("H85").Value = 2227.2727
$ws.Range("I85").Value = 1750
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 1750
$ws.Range("L85").Value = 2800
$ws.Range("M85").Value = -502
$ws.Range("N85").Value = -5296
$ws.Range("H88").Value = 28000
$ws.Range("J88").Value = 28000
$ws.Range("L88").Value = 28000
$ws.Range("N88").Value = -28856
$ws.Range("H91").Value = 28000
$ws.Range("J91").Value = 28000
$ws.Range("L91").Value = 28000
$ws.Range("N91").Value = -30964
$ws.Range("H122").Value = 4789.05
$ws.Range("I122").Value = 4710.9165
$ws.Range("J122").Value = 4906.25
$ws.Range("K122").Value = 14132.7495
$ws.Range("L122").Value = 14718.75
$ws.Range("M122").Value = -11682.7495
$ws.Range("N122").Value = -19618.75
$ws.Range("H126").Value = 8239.111000000001
$ws.Range("I126").Value = 21001.334
$ws.Range("J126").Value = 5686.6665
$ws.Range("K126").Value = 63004.00199999999
$ws.Range("L126").Value = 17059.9995
$ws.Range("M126").Value = -60534.00199999999
$ws.Range("N126").Value = -21999.9995
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 835.75
$ws.Range("I107").Value = 903.82355
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 2711.47065
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = -791.4706499999998
$ws.Range("N107").Value = -5190
$ws.Range("H122").Value = 1955.0333
$ws.Range("I122").Value = 2056.0833
$ws.Range("J122").Value = 1550.8334
$ws.Range("K122").Value = 6168.249899999999
$ws.Range("L122").Value = 4652.5002
$ws.Range("M122").Value = -3718.249899999999
$ws.Range("N122").Value = -9552.5002
$ws.Range("H126").Value = 4307
$ws.Range("I126").Value = 2867.9
$ws.Range("J126").Value = 11502.5
$ws.Range("K126").Value = 8603.700000000001
$ws.Range("L126").Value = 34507.5
$ws.Range("M126").Value = -6133.700000000001
$ws.Range("N126").Value = -39447.5
$ws.Range("H132").Value = 1318.4929
$ws.Range("I132").Value = 1044.6721
$ws.Range("J132").Value = 2988.8
$ws.Range("K132").Value = 3134.0163
$ws.Range("L132").Value = 8966.400000000001
$ws.Range("M132").Value = -604.0163000000002
$ws.Range("N132").Value = -14026.4
